$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume(1h) (E) columns with latest scraped values.
# D-column cells are forced to Text format while assigning so that Excel
# does not reinterpret numeric-looking strings (e.g. "246.46", "84.00",
# "0.06394") as actual numbers; the style is reset to Normal afterwards so
# no visible formatting/style change is introduced.

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.887.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +8.04%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.817.34"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +5.36%  "
$ws.Range("E4").Value = "  +0.16%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "246.46"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.56%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.12%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4934"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.44%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2776"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +7.22%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06394"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.49%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.810.75"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +4.97%  "
$ws.Range("E11").Value = "  +5.93%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.07068"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.6431"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +6.75%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "84.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +9.34%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.671"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +4.81%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "28.911.17"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +8.85%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.9997"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.13%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000007304"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.37%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.9992"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.12%  "
$ws.Range("E20").Value = "  +7.60%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "2.041.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +4.95%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "4.552"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +3.25%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "8.800"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +3.44%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "5.353"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +5.97%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.86"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.66%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "129.21"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +21.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.38"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.72%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "1.880"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +6.32%  "
$ws.Range("E29").Value = "  +2.85%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "4.124"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.24%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08343"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +5.35%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.779"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +3.15%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.04910"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +8.52%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.096"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +9.92%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.694"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +3.91%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.6728"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +9.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.291"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +15.11%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.701"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +10.29%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.9464"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +1.98%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.156"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +9.84%  "
$ws.Range("E41").Value = "  +5.83%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.9996"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +0.18%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "100.53"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.73%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.4084"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.79%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "7.163"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.86%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.1223"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +6.07%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.05515"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "8.091"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +2.32%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "31.62"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +5.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "1.304"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.80%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.3613"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +8.11%  "
